# Update testTakeABreak.xlsx: update size of takeABreakProgress
#
# 1) conditionTrials value changes from 3 to 4 (both columns B and C).
# 2) A new boolean parameter "responseTypedBool" (TRUE/TRUE) is inserted,
#    in its correct alphabetical position right after "conditionTrials".
# 3) "showTakeABreakCreditBool" column B changes from TRUE to FALSE
#    (column C stays TRUE).
# 4) A new numeric parameter "takeABreakMinimumDurationSec" (2/2) is
#    inserted, in its correct alphabetical position right before
#    "takeABreakTrialCredit".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- conditionTrials: 3 -> 4 -----------------------------------------
$ws.Range("B7").Value = 4
$ws.Range("C7").Value = 4

# --- insert "responseTypedBool" as the new row 8 ---------------------
$ws.Rows(8).Insert()
$ws.Range("A8").Value = "responseTypedBool"
$ws.Range("B8").Value = $true
$ws.Range("C8").Value = $true

# --- showTakeABreakCreditBool now lives on row 11; tweak column B ----
$ws.Range("B11").Value = $false

# --- insert "takeABreakMinimumDurationSec" as the new row 14 ---------
$ws.Rows(14).Insert()
$ws.Range("A14").Value = "takeABreakMinimumDurationSec"
$ws.Range("B14").Value = 2
$ws.Range("C14").Value = 2

# --- keep the active selection in sync with the authored workbook ----
[void]$ws.Range("C12").Select()
